# Auto-generated Excel COM-interop script applying the Omega_Profits value updates.
# Each block targets one worksheet/row; values are set directly per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 11989.909
$ws.Range("I43").Value = 15399.167
$ws.Range("K43").Value = 15399.167
$ws.Range("M43").Value = -15330.167

# Row 70
$ws.Range("H70").Value = 51819.285
$ws.Range("J70").Value = 70969
$ws.Range("L70").Value = 212907
$ws.Range("N70").Value = -213447

# Row 73
$ws.Range("H73").Value = 51819.285
$ws.Range("J73").Value = 70969
$ws.Range("L73").Value = 212907
$ws.Range("N73").Value = -214779

# Row 88
$ws.Range("H88").Value = 4256.7646
$ws.Range("I88").Value = 4588.6
$ws.Range("J88").Value = 4118.5
$ws.Range("K88").Value = 4588.6
$ws.Range("L88").Value = 4118.5
$ws.Range("M88").Value = -4182.6
$ws.Range("N88").Value = -4930.5

# Row 91
$ws.Range("H91").Value = 4256.7646
$ws.Range("I91").Value = 4588.6
$ws.Range("J91").Value = 4118.5
$ws.Range("K91").Value = 4588.6
$ws.Range("L91").Value = 4118.5
$ws.Range("M91").Value = -3184.6
$ws.Range("N91").Value = -6926.5

# Row 112
$ws.Range("H112").Value = 5876.375
$ws.Range("J112").Value = 6501.5713
$ws.Range("L112").Value = 19504.7139
$ws.Range("N112").Value = -21720.7139

# Row 125
$ws.Range("H125").Value = 900.4
$ws.Range("I125").Value = 929.125
$ws.Range("J125").Value = 785.5
$ws.Range("K125").Value = 8362.125
$ws.Range("L125").Value = 7069.5
$ws.Range("M125").Value = -5902.125
$ws.Range("N125").Value = -11989.5

# Row 138
$ws.Range("H138").Value = 3373.4775
$ws.Range("J138").Value = 4854.3
$ws.Range("L138").Value = 14562.9
$ws.Range("N138").Value = -24842.9

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1908.0869
$ws.Range("I2").Value = 1268.1111
$ws.Range("K2").Value = 1268.1111
$ws.Range("M2").Value = -1155.1111

# Row 33
$ws.Range("H33").Value = 48891.668
$ws.Range("I33").Value = 48891.668
$ws.Range("K33").Value = 48891.668
$ws.Range("M33").Value = -48562.668

# Row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = ""

# Row 61
$ws.Range("H61").Value = 5663.9546
$ws.Range("I61").Value = 5057.4375
$ws.Range("K61").Value = 5057.4375
$ws.Range("M61").Value = -4845.4375

# Row 74
$ws.Range("H74").Value = 1407.0571
$ws.Range("I74").Value = 1383.36
$ws.Range("J74").Value = 1466.3
$ws.Range("K74").Value = 1383.36
$ws.Range("L74").Value = 1466.3
$ws.Range("M74").Value = -509.3599999999999
$ws.Range("N74").Value = -3214.3

# Row 77
$ws.Range("H77").Value = 1407.0571
$ws.Range("I77").Value = 1383.36
$ws.Range("J77").Value = 1466.3
$ws.Range("K77").Value = 6916.799999999999
$ws.Range("L77").Value = 7331.5
$ws.Range("M77").Value = -2548.799999999999
$ws.Range("N77").Value = -16067.5

# Row 102
$ws.Range("H102").Value = 2541.8667
$ws.Range("I102").Value = 2682.3076
$ws.Range("K102").Value = 2682.3076
$ws.Range("M102").Value = -1060.3076

# Row 110
$ws.Range("H110").Value = 1802.4166
$ws.Range("I110").Value = 1696.3
$ws.Range("K110").Value = 1696.3
$ws.Range("M110").Value = 348.7

# Row 116
$ws.Range("H116").Value = 1908.0869
$ws.Range("I116").Value = 1268.1111
$ws.Range("K116").Value = 1268.1111
$ws.Range("M116").Value = 1025.8889

# Row 122
$ws.Range("H122").Value = 3888.625
$ws.Range("I122").Value = 3849
$ws.Range("J122").Value = 3954.6667
$ws.Range("K122").Value = 11547
$ws.Range("L122").Value = 11864.0001
$ws.Range("M122").Value = -9097
$ws.Range("N122").Value = -16764.0001

# Row 132
$ws.Range("H132").Value = 3340.4614
$ws.Range("I132").Value = 2691.65
$ws.Range("J132").Value = 5503.1665
$ws.Range("K132").Value = 8074.950000000001
$ws.Range("L132").Value = 16509.4995
$ws.Range("M132").Value = -5544.950000000001
$ws.Range("N132").Value = -21569.4995

# Row 136
$ws.Range("H136").Value = 5663.9546
$ws.Range("I136").Value = 5057.4375
$ws.Range("K136").Value = 15172.3125
$ws.Range("M136").Value = -12622.3125

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1908.0869
$ws.Range("I3").Value = 1268.1111
$ws.Range("K3").Value = 1268.1111
$ws.Range("M3").Value = -1154.1111

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 575
$ws.Range("I22").Value = 575
$ws.Range("K22").Value = 575
$ws.Range("M22").Value = -225

# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = ""

# Row 58
$ws.Range("H58").Value = 1772.5385
$ws.Range("I58").Value = 1943.8
$ws.Range("J58").Value = 1201.6666
$ws.Range("K58").Value = 1943.8
$ws.Range("L58").Value = 1201.6666
$ws.Range("M58").Value = -1740.8
$ws.Range("N58").Value = -1607.6666

# Row 107
$ws.Range("H107").Value = 2747.2
$ws.Range("I107").Value = 2904.2307
$ws.Range("K107").Value = 2904.2307
$ws.Range("M107").Value = -984.2307000000001

# Row 132
$ws.Range("H132").Value = 5630.9165
$ws.Range("I132").Value = 5086.3335
$ws.Range("K132").Value = 15259.0005
$ws.Range("M132").Value = -12729.0005

# Row 136
$ws.Range("H136").Value = 1772.5385
$ws.Range("I136").Value = 1943.8
$ws.Range("J136").Value = 1201.6666
$ws.Range("K136").Value = 5831.4
$ws.Range("L136").Value = 3604.9998
$ws.Range("M136").Value = -3281.4
$ws.Range("N136").Value = -8704.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 143.33333
$ws.Range("J11").Value = 500
$ws.Range("L11").Value = 1500
$ws.Range("N11").Value = -1780

# Row 75
$ws.Range("H75").Value = 1862
$ws.Range("I75").Value = 2479.3333
$ws.Range("J75").Value = 1719.5385
$ws.Range("K75").Value = 7437.999899999999
$ws.Range("L75").Value = 5158.6155
$ws.Range("M75").Value = -6439.999899999999
$ws.Range("N75").Value = -7154.6155

# Row 78
$ws.Range("H78").Value = 1862
$ws.Range("I78").Value = 2479.3333
$ws.Range("J78").Value = 1719.5385
$ws.Range("K78").Value = 22313.9997
$ws.Range("L78").Value = 15475.8465
$ws.Range("M78").Value = -17321.9997
$ws.Range("N78").Value = -25459.8465

# Row 92
$ws.Range("H92").Value = 366.16666
$ws.Range("I92").Value = 299.66666
$ws.Range("J92").Value = 432.66666
$ws.Range("K92").Value = 898.9999799999999
$ws.Range("L92").Value = 1297.99998
$ws.Range("M92").Value = 349.0000200000001
$ws.Range("N92").Value = -3793.99998

# Row 98
$ws.Range("H98").Value = 1640.5
$ws.Range("I98").Value = 1189.8
$ws.Range("K98").Value = 3569.4
$ws.Range("M98").Value = -2071.4

# Row 131
$ws.Range("H131").Value = 2043.762
$ws.Range("I131").Value = 1245.3846
$ws.Range("K131").Value = 3736.1538
$ws.Range("M131").Value = 1303.8462

$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 173616.17
$ws.Range("I14").Value = 250799.5
$ws.Range("K14").Value = 250799.5
$ws.Range("M14").Value = -250631.5

# Row 122
$ws.Range("H122").Value = 1865.7142
$ws.Range("I122").Value = 1865.7142
$ws.Range("K122").Value = 5597.142599999999
$ws.Range("M122").Value = -3147.142599999999

$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 250584.25
$ws.Range("I14").Value = 250584.25
$ws.Range("K14").Value = 250584.25
$ws.Range("M14").Value = -250412.25

# Row 16
$ws.Range("H16").Value = 1252.8235
$ws.Range("I16").Value = 1018.625
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 1018.625
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -848.625
$ws.Range("N16").Value = -5340

# Row 55
$ws.Range("H55").Value = 278.29413
$ws.Range("I55").Value = 344.44446
$ws.Range("J55").Value = 203.875
$ws.Range("K55").Value = 344.44446
$ws.Range("L55").Value = 203.875
$ws.Range("M55").Value = -171.44446
$ws.Range("N55").Value = -549.875

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 3778.5625
$ws.Range("I107").Value = 5376
$ws.Range("J107").Value = 2181.125
$ws.Range("K107").Value = 16128
$ws.Range("L107").Value = 6543.375
$ws.Range("M107").Value = -14208
$ws.Range("N107").Value = -10383.375

# Row 122
$ws.Range("H122").Value = 2653.8572
$ws.Range("I122").Value = 2740.0667
$ws.Range("K122").Value = 8220.2001
$ws.Range("M122").Value = -5770.2001

# Row 126
$ws.Range("H126").Value = 2465.8333
$ws.Range("I126").Value = 1509.4
$ws.Range("K126").Value = 4528.200000000001
$ws.Range("M126").Value = -2058.200000000001

# Row 136
$ws.Range("H136").Value = 2386.9666
$ws.Range("I136").Value = 1768.64
$ws.Range("K136").Value = 5305.92
$ws.Range("M136").Value = -2755.92

